$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 5.5
$ws.Range("AA2").Value = 13
$ws.Range("AC2").Value = 10
$ws.Range("L3").Value = 2.5
$ws.Range("O3").Value = 1.29
$ws.Range("P3").Value = 3.75
$ws.Range("Q3").Value = 1.91
$ws.Range("R3").Value = 1.99
$ws.Range("W3").Value = 1.36
$ws.Range("X3").Value = 3
$ws.Range("Y3").Value = 1.7
$ws.Range("Z3").Value = 2.05
$ws.Range("AA3").Value = 13
$ws.Range("AQ3").Value = 23
$ws.Range("G5").Value = 2.15
$ws.Range("I5").Value = 3.9
$ws.Range("J5").Value = 2.88
$ws.Range("L5").Value = 4.5
$ws.Range("AD5").Value = 19
$ws.Range("Q9").Value = 2.5
$ws.Range("R9").Value = 1.5
$ws.Range("U9").Value = 5
$ws.Range("V9").Value = 1.17
$ws.Range("AR9").Value = 1.95
$ws.Range("AS9").Value = 1.9
$ws.Range("H12").Value = 7.9
$ws.Range("I12").Value = 25
$ws.Range("J12").Value = 1.4
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 15
$ws.Range("Q12").Value = 1.42
$ws.Range("R12").Value = 2.7
$ws.Range("W12").Value = 1.24
$ws.Range("X12").Value = 3.75
$ws.Range("Y12").Value = 2.62
$ws.Range("Z12").Value = 1.44
$ws.Range("AA12").Value = 6.8
$ws.Range("AB12").Value = 6
$ws.Range("AJ12").Value = 350
$ws.Range("AL12").Value = 40
$ws.Range("AM12").Value = 350
$ws.Range("AN12").Value = 100
$ws.Range("AP12").Value = 700
$ws.Range("O14").Value = 1.22
$ws.Range("P14").Value = 4
$ws.Range("U14").Value = 2.75
$ws.Range("V14").Value = 1.4
$ws.Range("G15").Value = 2.45
$ws.Range("I15").Value = 2.88
$ws.Range("J15").Value = 3.1
$ws.Range("O15").Value = 1.3
$ws.Range("P15").Value = 3.5
$ws.Range("Q15").Value = 1.98
$ws.Range("R15").Value = 1.88
$ws.Range("Y15").Value = 1.73
$ws.Range("Z15").Value = 2
$ws.Range("AF15").Value = 29
$ws.Range("AH15").Value = 6
$ws.Range("AL15").Value = 9.5
$ws.Range("G16").Value = 1.65
$ws.Range("J16").Value = 2.3
$ws.Range("Q16").Value = 2.35
$ws.Range("R16").Value = 1.57
$ws.Range("U16").Value = 4.5
$ws.Range("V16").Value = 1.18
$ws.Range("AG16").Value = 7.5
$ws.Range("AH16").Value = 7.5
$ws.Range("AJ16").Value = 101
$ws.Range("AR16").Value = 1.8
$ws.Range("AS16").Value = 2.05
$ws.Range("G17").Value = 2.6
$ws.Range("I17").Value = 2.3
$ws.Range("J17").Value = 3.2
$ws.Range("L17").Value = 2.88
$ws.Range("AA17").Value = 13
$ws.Range("AB17").Value = 17
$ws.Range("AC17").Value = 11
$ws.Range("AD17").Value = 29
$ws.Range("AF17").Value = 23
$ws.Range("AH17").Value = 7
$ws.Range("AO17").Value = 23
$ws.Range("AP17").Value = 17
$ws.Range("G19").Value = 3.9
$ws.Range("H19").Value = 3.6
$ws.Range("I19").Value = 1.8
$ws.Range("J19").Value = 4.15
$ws.Range("K19").Value = 2.22
$ws.Range("L19").Value = 2.32
$ws.Range("O19").Value = 1.21
$ws.Range("P19").Value = 3.55
$ws.Range("Q19").Value = 1.65
$ws.Range("R19").Value = 2
$ws.Range("U19").Value = 2.47
$ws.Range("V19").Value = 1.42
$ws.Range("Y19").Value = 1.6
$ws.Range("Z19").Value = 2.07
$ws.Range("AA19").Value = 13.5
$ws.Range("AB19").Value = 23
$ws.Range("AC19").Value = 13
$ws.Range("AD19").Value = 60
$ws.Range("AF19").Value = 35
$ws.Range("AG19").Value = 12.5
$ws.Range("AH19").Value = 7.2
$ws.Range("AJ19").Value = 50
$ws.Range("AL19").Value = 8.5
$ws.Range("AO19").Value = 15.5
$ws.Range("AP19").Value = 13
$ws.Range("AQ19").Value = 21
$ws.Range("G20").Value = 1.65
$ws.Range("H20").Value = 3.9
$ws.Range("I20").Value = 4.75
$ws.Range("J20").Value = 2.2
$ws.Range("N20").Value = 17
$ws.Range("S20").Value = 1.93
$ws.Range("T20").Value = 1.93
$ws.Range("W20").Value = 1.29
$ws.Range("X20").Value = 3.5
$ws.Range("Y20").Value = 1.57
$ws.Range("Z20").Value = 2.25
$ws.Range("AD20").Value = 13
$ws.Range("AE20").Value = 12
$ws.Range("AG20").Value = 17
$ws.Range("AH20").Value = 8
$ws.Range("AM20").Value = 29
$ws.Range("G21").Value = 2.2
$ws.Range("H21").Value = 3.6
$ws.Range("L21").Value = 3.4
$ws.Range("O21").Value = 1.18
$ws.Range("P21").Value = 4.5
$ws.Range("Q21").Value = 1.62
$ws.Range("R21").Value = 2.25
$ws.Range("U21").Value = 2.5
$ws.Range("V21").Value = 1.5
$ws.Range("AB21").Value = 13
$ws.Range("G22").Value = 1.95
$ws.Range("H22").Value = 3.8
$ws.Range("J22").Value = 2.5
$ws.Range("K22").Value = 2.4
$ws.Range("N22").Value = 17
$ws.Range("Q22").Value = 1.53
$ws.Range("R22").Value = 2.4
$ws.Range("S22").Value = 1.88
$ws.Range("T22").Value = 1.98
$ws.Range("U22").Value = 2.25
$ws.Range("V22").Value = 1.57
$ws.Range("Y22").Value = 1.5
$ws.Range("Z22").Value = 2.5
$ws.Range("AK22").Value = 101
$ws.Range("N23").Value = 17
$ws.Range("S23").Value = 1.9
$ws.Range("T23").Value = 1.95
$ws.Range("U23").Value = 2.25
$ws.Range("V23").Value = 1.57
$ws.Range("G25").Value = 3.7
$ws.Range("H25").Value = 3.3
$ws.Range("I25").Value = 2.05
$ws.Range("K25").Value = 2.05
$ws.Range("M25").Value = 1.07
$ws.Range("N25").Value = 8.5
$ws.Range("Q25").Value = 2.1
$ws.Range("R25").Value = 1.7
$ws.Range("AG25").Value = 8.5
$ws.Range("M26").Value = 1.03
$ws.Range("N26").Value = 15
$ws.Range("Q26").Value = 1.7
$ws.Range("R26").Value = 2.1
$ws.Range("G27").Value = 1.82
$ws.Range("I27").Value = 4
$ws.Range("J27").Value = 2.42
$ws.Range("K27").Value = 2.15
$ws.Range("L27").Value = 4.25
$ws.Range("AA27").Value = 7.6
$ws.Range("AD27").Value = 15
$ws.Range("AL27").Value = 13
$ws.Range("AM27").Value = 24
$ws.Range("AN27").Value = 13
$ws.Range("AO27").Value = 60
$ws.Range("AP27").Value = 35
$ws.Range("G28").Value = 1.5
$ws.Range("I28").Value = 7
$ws.Range("AB28").Value = 6.5
$ws.Range("AD28").Value = 10
$ws.Range("AN28").Value = 21
$ws.Range("AG30").Value = 7.3
